$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Snapshot the existing 33 data rows (rows 2..34) keyed by the id in col A.
#    We read through COM so we get the real resolved text (not raw XML), which
#    keeps any literal "&amp;" / unicode / mojibake sequences intact exactly
#    as Excel exposes them. Two parallel hashtables (tweet/label) avoid any
#    ambiguity around nested-array assignment.
# ---------------------------------------------------------------------------
$tweetById = @{}
$labelById = @{}
for ($r = 2; $r -le 34; $r++) {
    $id = $ws.Cells.Item($r, 1).Value2
    $tweetById[[int]$id] = $ws.Cells.Item($r, 2).Value2
    $labelById[[int]$id] = $ws.Cells.Item($r, 3).Value2
}

# ---------------------------------------------------------------------------
# 2. New row order (by id), row 1 .. row 33. This is the same 33 records,
#    just re-sequenced; the former header row (id/tweet/label) is gone and
#    row 34 becomes a blank trailing row.
# ---------------------------------------------------------------------------
$newOrder = @(6392,6396,6393,6391,6394,6390,6395,6389,6397,6388,6398,6399,6400,6401,6407,6402,6408,6403,6409,6404,6410,6405,6411,6406,6413,6412,6414,6415,6416,6417,6418,6419,6420)

for ($i = 0; $i -lt $newOrder.Count; $i++) {
    $r = $i + 1
    $id = $newOrder[$i]
    $ws.Cells.Item($r, 1).Value = $id
    $ws.Cells.Item($r, 2).Value = $tweetById[[int]$id]
    $ws.Cells.Item($r, 3).Value = $labelById[[int]$id]
}

# ---------------------------------------------------------------------------
# 3. Row 34 becomes an empty row (values cleared, formatting/style kept).
# ---------------------------------------------------------------------------
$ws.Range("A34:C34").ClearContents()

# ---------------------------------------------------------------------------
# 4. Row height tweaks on rows 27 and 31 (14.25pt custom height).
# ---------------------------------------------------------------------------
$ws.Rows.Item(27).RowHeight = 14.25
$ws.Rows.Item(31).RowHeight = 14.25

# ---------------------------------------------------------------------------
# 5. Sheet view: scroll so column B is left-most visible, and move the active
#    selection to B8.
# ---------------------------------------------------------------------------
$ws.Range("B1").Select()
$excel.ActiveWindow.ScrollColumn = 2
$ws.Range("B8").Select()
